$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting of row 15 down into the new row 16
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in new row 16 content
$ws.Range("A16").Value = "tap15"
$ws.Range("B16").Value = 0.0033
$ws.Range("C16").Formula = "=2^15"
$ws.Range("D16").Formula = "=B16*C16"
$ws.Range("E16").Formula = "=D16"
$ws.Range("F16").Formula = "=E16"
$ws.Range("G16").Formula = "=DEC2HEX(F16)"

# Update existing coefficients (column B) with new values
$ws.Range("B1").Value = 0.0033
$ws.Range("B2").Value = -0.0136
$ws.Range("B3").Value = -0.0053
$ws.Range("B4").Value = 0.0408
$ws.Range("B5").Value = -0.0104
$ws.Range("B6").Value = -0.1022
$ws.Range("B7").Value = 0.0992
$ws.Range("B8").Value = 0.4854
$ws.Range("B9").Value = 0.4854
$ws.Range("B10").Value = 0.0992
$ws.Range("B11").Value = -0.1022
$ws.Range("B12").Value = -0.0104
$ws.Range("B13").Value = 0.0408
$ws.Range("B14").Value = -0.0053
$ws.Range("B15").Value = -0.0136

# Update the selection to reflect the new active cell
$ws.Range("B16").Select() | Out-Null
